$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the cryptos.xlsx diff (89 cell text changes).
# Column D "Price" values are plain text that often LOOK numeric (e.g. "6.00",
# "12.70", "0.0000210"). A straight .Value assignment lets Excel auto-convert
# those into real numbers and silently drop the formatting-significant trailing
# zeros, which would diverge from the source workbook (inline strings). For any
# such cell we briefly force a Text number format, assign the literal string,
# then restore the cell to the "Normal" style so no stray formatting is left
# behind (matches the target XML, which carries no style attribute on these cells).

$ws.Range("D2").Value = "57.443.43"
$ws.Range("E2").Value = "  -7.12%  "
$ws.Range("D3").Value = "2.886.70"
$ws.Range("E3").Value = "  -5.51%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "121.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.91%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "2.878.13"
$ws.Range("E8").Value = "  -5.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.488"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.124"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -10.82%  "
$ws.Range("E11").Value = "  -9.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.429"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000210"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -11.31%  "
$ws.Range("E14").Value = "  -7.56%  "
$ws.Range("E15").Value = "  -0.65%  "
$ws.Range("D16").Value = "3.362.14"
$ws.Range("E16").Value = "  -5.58%  "
$ws.Range("D17").Value = "2.881.33"
$ws.Range("E17").Value = "  -5.86%  "
$ws.Range("D18").Value = "57.346.60"
$ws.Range("E18").Value = "  -7.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "405.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -9.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.648"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.87%  "
$ws.Range("E23").Value = "  -8.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "76.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.71%  "
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("E28").Value = "  -5.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "24.49"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0946"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.21%  "
$ws.Range("E34").Value = "  -13.75%  "
$ws.Range("B35").Value = "Mantle"
$ws.Range("C35").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.894"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.48%  "
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.33"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "48.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.28%  "
$ws.Range("E38").Value = "  +4.54%  "
$ws.Range("D39").Value = "0.0₃0612"
$ws.Range("E39").Value = "  -12.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0341"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.28%  "
$ws.Range("E41").Value = "  -3.43%  "
$ws.Range("D42").Value = "2.607.51"
$ws.Range("E42").Value = "  -3.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "355.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.43%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "117.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.69%  "
$ws.Range("E47").Value = "  -6.04%  "
$ws.Range("E48").Value = "  -3.52%  "
$ws.Range("E49").Value = "  -5.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.25"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.52%  "
